$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Oval 3" shape (the decorative ellipse on the title slide).
$orig = $s.Shapes.Item(3)

# Clear its (empty) text so the stray paragraph-mark formatting
# (line/fill/shadow) carried on the endParaRPr is dropped.
$orig.TextFrame.TextRange.Text = ""

# Duplicate it - PowerPoint hands the copy a fresh shape id (5) while the
# original keeps id 4; deleting the original afterwards leaves the new
# shape as the sole survivor with id 5, matching the recorded edit.
$dupColl = $orig.Duplicate()
$dup = $dupColl.Item(1)
$orig.Delete()

# Rename + reshape: ellipse -> rectangle.
$dup.Name = "Rectangle 4"
$dup.AutoShapeType = 1

# Reposition / resize to match the new layout (values are the exact
# point equivalents of the target EMU offsets/extents).
$dup.Left = 44.44448818897638
$dup.Top = 280.95826771653543
$dup.Width = 199.1111073622047
$dup.Height = 170.6666949133858
